# Prepping heat model output data: summarize rasterization data and
# calculate added heat from vehicles and pavement across metro.
#
# - Corrects the "car SUV" UDDS CS share in C29 (0.41 -> 0.42), which
#   also updates the SUM() total in H29 from 0.99 to 1 automatically.
# - Adds a small summary table (rows 34-35) with the drivecycle headers
#   (UDDS CS, HWY, US06, SC03) and the average mpg across all vehicle
#   classes for each drivecycle.
# - Updates the window scroll position / selection to reflect the newly
#   added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the car SUV UDDS CS share.
$ws.Range("C29").Value = 0.42

# New drivecycle header row (reuses existing shared strings: UDDS CS,
# HWY, US06, SC03).
$ws.Range("B34").Value = "UDDS CS"
$ws.Range("C34").Value = "HWY"
$ws.Range("D34").Value = "US06"
$ws.Range("E34").Value = "SC03"

# New row with the average mpg per drivecycle, across the four vehicle
# classes (rows 2/14/18/22 for UDDS CS, etc.).
$ws.Range("B35").Formula = "=AVERAGE(C2,C14,C18,C22)"
$ws.Range("C35").Formula = "=AVERAGE(C3,C15,C19,C23)"
$ws.Range("D35").Formula = "=AVERAGE(C4,C16,C20,C24)"
$ws.Range("E35").Formula = "=AVERAGE(C5,C17,C21,C25)"

# Scroll the view down to the new data and select the cell just past it.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E36").Select()
